$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date serial 45181 (2023-09-12) for every
# data row from row 2 to row 77. Update it to 45182 (2023-09-13).
for ($row = 2; $row -le 77; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value = 45182
    }
}
